# Update the appointment-time values in the Appointments sheet.
# Row 2 (Appointment ID 1): normalise the "pm" suffix to "PM".
# Row 3 (Appointment ID 2): the appointment was rescheduled from
#   18-Nov-2024 4:00 pm to 20-Nov-2024 5:00 PM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "18-Nov-2024 2:00:00 PM"
$ws.Range("D3").Value = "20-Nov-2024 5:00:00 PM"

# Leave the selection on the cell that was last edited, matching where
# the workbook was left when it was saved.
$ws.Range("D2").Select()
